# "Final Work Day Commit!"
#
# Mark a further batch of pages as checked in the export tracker.
# Column C = "Success" marker, column E = "Problem" marker; both use the
# literal text "x". Rows 89-124 (inclusive) each gain exactly one "x" in
# either C or E, matching the day's QA pass. The COUNTIF-driven summary
# cells D1/F1 recalculate automatically once the new marks are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E89").Value = "x"
$ws.Range("E90").Value = "x"
$ws.Range("C91").Value = "x"
$ws.Range("E92").Value = "x"
$ws.Range("C93").Value = "x"
$ws.Range("E94").Value = "x"
$ws.Range("E95").Value = "x"
$ws.Range("C96").Value = "x"
$ws.Range("C97").Value = "x"
$ws.Range("E98").Value = "x"
$ws.Range("C99").Value = "x"
$ws.Range("C100").Value = "x"
$ws.Range("C101").Value = "x"
$ws.Range("C102").Value = "x"
$ws.Range("C103").Value = "x"
$ws.Range("C104").Value = "x"
$ws.Range("C105").Value = "x"
$ws.Range("C106").Value = "x"
$ws.Range("C107").Value = "x"
$ws.Range("C108").Value = "x"
$ws.Range("C109").Value = "x"
$ws.Range("E110").Value = "x"
$ws.Range("E111").Value = "x"
$ws.Range("C112").Value = "x"
$ws.Range("C113").Value = "x"
$ws.Range("C114").Value = "x"
$ws.Range("C115").Value = "x"
$ws.Range("C116").Value = "x"
$ws.Range("C117").Value = "x"
$ws.Range("C118").Value = "x"
$ws.Range("C119").Value = "x"
$ws.Range("E120").Value = "x"
$ws.Range("E121").Value = "x"
$ws.Range("C122").Value = "x"
$ws.Range("C123").Value = "x"
$ws.Range("E124").Value = "x"

# Leave the selection where the author stopped for the day (the sheet's
# frozen header pane is already in place from the prior session).
$null = $ws.Range("B125").Select()
